$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4 = 0.353
    5 = 0.712
    6 = 0.868
    7 = 0.806
    8 = 0.679
    9 = 0.592
    10 = 0.421
    11 = 0.199
    27 = 0.056
    28 = 0.263
    29 = 0.446
    30 = 0.653
    31 = 0.724
    32 = 0.714
    33 = 0.64
    34 = 0.395
    36 = 0.026
    51 = 0.041
    52 = 0.173
    54 = 0.491
    55 = 0.592
    56 = 0.526
    57 = 0.273
    58 = 0.172
    59 = 0.063
    60 = 0.012
    76 = 0.347
    77 = 0.582
    79 = 1.063
    80 = 1.006
    81 = 0.8179999999999999
    82 = 0.504
    83 = 0.172
    84 = 0.025
    99 = 0.195
    102 = 2.388
    104 = 2.524
    105 = 1.921
    106 = 1.151
    148 = 0.968
    149 = 1.986
    150 = 2.805
    151 = 2.985
    152 = 2.73
    153 = 2.215
    154 = 1.127
    155 = 0.375
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item([int]$row, 3).Value = $updates[$row]
}
